$d = $word.ActiveDocument
$d.Content.Find.Execute("showing 6 pips, then 1 pip", $true, $false, $false, $false, $false,
                         $true, 1, $false, "showing 4 pips, then 1 pip", 2)
